# Update data and freezer boxes
# Target sheets: "R4P7" (rId9 old active tab) and "R4_P9" (new active tab)

$wb = $excel.ActiveWorkbook
$wsP7 = $wb.Worksheets.Item("R4P7")
$wsP9 = $wb.Worksheets.Item("R4_P9")

# ------------------------------------------------------------------
# R4_P9 ("Sheet17") data-table updates
# ------------------------------------------------------------------

# Row 9 continues the A2-13-x series into columns I/J
$wsP9.Cells.Item(9, 9).Value  = "A2-13-6"
$wsP9.Cells.Item(9, 10).Value = "A2-13-7"

# Row 10: A2-13-8 .. A2-13-15
$row10 = @("A2-13-8","A2-13-9","A2-13-10","A2-13-11","A2-13-12","A2-13-13","A2-13-14","A2-13-15")
$col = 3
foreach ($v in $row10) {
    $wsP9.Cells.Item(10, $col).Value = $v
    $col++
}

# Row 11: A2-14-1 .. A2-14-8
$row11 = @("A2-14-1","A2-14-2","A2-14-3","A2-14-4","A2-14-5","A2-14-6","A2-14-7","A2-14-8")
$col = 3
foreach ($v in $row11) {
    $wsP9.Cells.Item(11, $col).Value = $v
    $col++
}

# Row 12: A2-14-9, then (out of sequence) A2-15-9 lands in I15, then A2-14-10 .. A2-14-16
$wsP9.Cells.Item(12, 3).Value = "A2-14-9"
$wsP9.Cells.Item(15, 9).Value = "A2-15-9"
$row12rest = @("A2-14-10","A2-14-11","A2-14-12","A2-14-13","A2-14-14","A2-14-15","A2-14-16")
$col = 4
foreach ($v in $row12rest) {
    $wsP9.Cells.Item(12, $col).Value = $v
    $col++
}

# Row 13: A2-14-17 .. A2-14-24
$row13 = @("A2-14-17","A2-14-18","A2-14-19","A2-14-20","A2-14-21","A2-14-22","A2-14-23","A2-14-24")
$col = 3
foreach ($v in $row13) {
    $wsP9.Cells.Item(13, $col).Value = $v
    $col++
}

# Row 14: A2-14-25 .. A2-14-30, A2-15-1, A2-15-2
$row14 = @("A2-14-25","A2-14-26","A2-14-27","A2-14-28","A2-14-29","A2-14-30","A2-15-1","A2-15-2")
$col = 3
foreach ($v in $row14) {
    $wsP9.Cells.Item(14, $col).Value = $v
    $col++
}

# Row 15: A2-15-3 .. A2-15-8 (C..H), I15 already holds A2-15-9, then J15 = A2-15-10
$row15 = @("A2-15-3","A2-15-4","A2-15-5","A2-15-6","A2-15-7","A2-15-8")
$col = 3
foreach ($v in $row15) {
    $wsP9.Cells.Item(15, $col).Value = $v
    $col++
}
$wsP9.Cells.Item(15, 10).Value = "A2-15-10"

# Row 16: A2-15-11 .. A2-15-18
$row16 = @("A2-15-11","A2-15-12","A2-15-13","A2-15-14","A2-15-15","A2-15-16","A2-15-17","A2-15-18")
$col = 3
foreach ($v in $row16) {
    $wsP9.Cells.Item(16, $col).Value = $v
    $col++
}

# ------------------------------------------------------------------
# R4_P9 note cells (B2 / B6) get amended text
# ------------------------------------------------------------------
$wsP9.Range("B2").Value = " cohort 2018-09-11(day 12+27) [1 sample]; 2018-09-12 cohort (day 13) ; 2018-09-13 cohor t(day 14); 2018-09-14 cohort (day 15)"
$wsP9.Range("B6").Value = "NB#005, pg 24,25, 44"

# ------------------------------------------------------------------
# R4_P9 row 21 note ("is in different box") is cleared out
# ------------------------------------------------------------------
$wsP9.Cells.Item(21, 2).ClearContents()
$wsP9.Cells.Item(21, 3).ClearContents()

# ------------------------------------------------------------------
# Active-tab / selection bookkeeping: selection moves from R4P7 to R4_P9
# ------------------------------------------------------------------
$wsP9.Activate()
$wsP9.Range("B7").Select()
